# Automatic update of files.
# Remove the superfluous "test" friendly-name argument from the HYPERLINK()
# formulas in row 2 (S2, T2, V2, W2, X2, Y2), leaving only the URL argument.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/artfynd/A 34293-2023.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/kartor/A 34293-2023.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/klagomål/A 34293-2023.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/klagomålsmail/A 34293-2023.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/tillsyn/A 34293-2023.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/tillsynsmail/A 34293-2023.docx")'
